# Auto-generated Excel COM-interop script to apply Chocobo_Profits market-data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 357.94116
$ws.Range("I5").Value = 231.66667
$ws.Range("K5").Value = 231.66667
$ws.Range("M5").Value = -116.66667
$ws.Range("H12").Value = 1381
$ws.Range("J12").Value = 114
$ws.Range("L12").Value = 114
$ws.Range("N12").Value = -454
$ws.Range("H38").Value = 3067.9473
$ws.Range("I38").Value = 162.625
$ws.Range("J38").Value = 5180.909
$ws.Range("K38").Value = 487.875
$ws.Range("L38").Value = 15542.727
$ws.Range("M38").Value = -115.875
$ws.Range("N38").Value = -16286.727
$ws.Range("H40").Value = 2275
$ws.Range("I40").Value = 2275
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2275
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -2100
$ws.Range("N40").ClearContents()
$ws.Range("H62").Value = 1243.2222
$ws.Range("I62").Value = 1369.8572
$ws.Range("J62").Value = 800
$ws.Range("K62").Value = 1369.8572
$ws.Range("L62").Value = 800
$ws.Range("M62").Value = -745.8571999999999
$ws.Range("N62").Value = -2048
$ws.Range("H64").Value = 2233.3333
$ws.Range("H65").Value = 1243.2222
$ws.Range("I65").Value = 1369.8572
$ws.Range("J65").Value = 800
$ws.Range("K65").Value = 6849.286
$ws.Range("L65").Value = 4000
$ws.Range("M65").Value = -3729.286
$ws.Range("N65").Value = -10240
$ws.Range("H67").Value = 2233.3333
$ws.Range("H88").Value = 21950
$ws.Range("J88").Value = 41400
$ws.Range("L88").Value = 41400
$ws.Range("N88").Value = -42212
$ws.Range("H91").Value = 21950
$ws.Range("J91").Value = 41400
$ws.Range("L91").Value = 41400
$ws.Range("N91").Value = -44208
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("H112").Value = 1317.5167
$ws.Range("J112").Value = 1331.3729
$ws.Range("L112").Value = 3994.1187
$ws.Range("N112").Value = -6210.1187
$ws.Range("H116").Value = 440893.97
$ws.Range("J116").Value = 8041.3125
$ws.Range("L116").Value = 8041.3125
$ws.Range("N116").Value = -14925.3125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 10200
$ws.Range("I35").Value = 1600
$ws.Range("J35").Value = 36000
$ws.Range("K35").Value = 1600
$ws.Range("L35").Value = 36000
$ws.Range("M35").Value = -1194
$ws.Range("N35").Value = -36812
$ws.Range("H45").Value = 3495.5
$ws.Range("I45").Value = 4006
$ws.Range("J45").Value = 1964
$ws.Range("K45").Value = 4006
$ws.Range("L45").Value = 1964
$ws.Range("M45").Value = -3629
$ws.Range("N45").Value = -2718

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2017.8379
$ws.Range("I99").Value = 1051.6
$ws.Range("J99").Value = 4030.8333
$ws.Range("K99").Value = 1051.6
$ws.Range("L99").Value = 4030.8333
$ws.Range("M99").Value = 446.4000000000001
$ws.Range("N99").Value = -7026.8333
$ws.Range("H105").Value = 1823.4429
$ws.Range("I105").Value = 1825.8806
$ws.Range("J105").Value = 1769
$ws.Range("K105").Value = 1825.8806
$ws.Range("L105").Value = 1769
$ws.Range("M105").Value = -78.88059999999996
$ws.Range("N105").Value = -5263

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5353.75
$ws.Range("I31").Value = 1566.9062
$ws.Range("J31").Value = 11412.7
$ws.Range("K31").Value = 1566.9062
$ws.Range("L31").Value = 11412.7
$ws.Range("M31").Value = -1271.9062
$ws.Range("N31").Value = -12002.7
$ws.Range("H34").Value = 5353.75
$ws.Range("I34").Value = 1566.9062
$ws.Range("J34").Value = 11412.7
$ws.Range("K34").Value = 1566.9062
$ws.Range("L34").Value = 11412.7
$ws.Range("M34").Value = -1364.9062
$ws.Range("N34").Value = -11816.7
$ws.Range("H137").Value = 48373.332
$ws.Range("J137").Value = 48373.332
$ws.Range("L137").Value = 48373.332
$ws.Range("N137").Value = -58573.332

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 5435366.5
$ws.Range("I113").Value = 601.8570999999999
$ws.Range("J113").Value = 13889445
$ws.Range("K113").Value = 1805.5713
$ws.Range("L113").Value = 41668335
$ws.Range("M113").Value = 364.4287000000002
$ws.Range("N113").Value = -41672675
$ws.Range("H122").Value = 2783.647
$ws.Range("J122").Value = 3095.3022
$ws.Range("L122").Value = 27857.7198
$ws.Range("N122").Value = -32757.7198
$ws.Range("H131").Value = 782.9484
$ws.Range("J131").Value = 825.7977
$ws.Range("L131").Value = 2477.3931
$ws.Range("N131").Value = -12557.3931
$ws.Range("H132").Value = 2052.4666
$ws.Range("I132").Value = 820.6429000000001
$ws.Range("J132").Value = 3130.3125
$ws.Range("K132").Value = 7385.7861
$ws.Range("L132").Value = 28172.8125
$ws.Range("M132").Value = -4855.7861
$ws.Range("N132").Value = -33232.8125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 29006
$ws.Range("J20").Value = 29006
$ws.Range("L20").Value = 29006
$ws.Range("N20").Value = -29496
$ws.Range("H27").Value = 23666.666
$ws.Range("J27").Value = 33000
$ws.Range("L27").Value = 33000
$ws.Range("N27").Value = -33332
$ws.Range("H96").Value = 29440.25
$ws.Range("J96").Value = 29440.25
$ws.Range("L96").Value = 29440.25
$ws.Range("N96").Value = -34932.25
$ws.Range("H140").Value = 38610.555
$ws.Range("J140").Value = 38610.555
$ws.Range("L140").Value = 38610.555
$ws.Range("N140").Value = -48970.555

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3586.8
$ws.Range("I7").Value = 1723.5834
$ws.Range("J7").Value = 5306.6924
$ws.Range("K7").Value = 1723.5834
$ws.Range("L7").Value = 5306.6924
$ws.Range("M7").Value = -1611.5834
$ws.Range("N7").Value = -5530.6924
$ws.Range("H126").Value = 3586.8
$ws.Range("I126").Value = 1723.5834
$ws.Range("J126").Value = 5306.6924
$ws.Range("K126").Value = 5170.7502
$ws.Range("L126").Value = 15920.0772
$ws.Range("M126").Value = -2700.7502
$ws.Range("N126").Value = -20860.0772
$ws.Range("H135").Value = 49416.668
$ws.Range("J135").Value = 49416.668
$ws.Range("L135").Value = 49416.668
$ws.Range("N135").Value = -59556.668
$ws.Range("H139").Value = 48702
$ws.Range("J139").Value = 48702
$ws.Range("L139").Value = 48702
$ws.Range("N139").Value = -58982
$ws.Range("H140").Value = 58541.95
$ws.Range("J140").Value = 58541.95
$ws.Range("L140").Value = 58541.95
$ws.Range("N140").Value = -68901.95
$ws.Range("H141").Value = 32019.5
$ws.Range("J141").Value = 32019.5
$ws.Range("L141").Value = 32019.5
$ws.Range("N141").Value = -42379.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 43936092
$ws.Range("I96").Value = 56139824
$ws.Range("J96").Value = 2649.8
$ws.Range("K96").Value = 56139824
$ws.Range("L96").Value = 2649.8
$ws.Range("M96").Value = -56138451
$ws.Range("N96").Value = -5395.8
$ws.Range("H138").Value = 41999.332
$ws.Range("J138").Value = 41999.332
$ws.Range("L138").Value = 41999.332
$ws.Range("N138").Value = -52279.332
$ws.Range("H139").Value = 43028.75
$ws.Range("J139").Value = 43028.75
$ws.Range("L139").Value = 43028.75
$ws.Range("N139").Value = -53308.75
$ws.Range("H140").Value = 28529.25
$ws.Range("J140").Value = 28529.25
$ws.Range("L140").Value = 28529.25
$ws.Range("N140").Value = -38889.25
$ws.Range("H141").Value = 37469
$ws.Range("J141").Value = 37469
$ws.Range("L141").Value = 37469
$ws.Range("N141").Value = -47829
